$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild rows 2-10 with the updated TPM NATMI ligand-receptor data.
# Sending clusters: FAPs, MuSCs, Resolving-Mac; Target clusters: ECs, FAPs, MuSCs
$data = @{}

$data[2] = @{
    "A" = "FAPs"
    "B" = "Tac1"
    "C" = "Tacr1"
    "D" = "ECs"
    "E" = 3
    "F" = 1
    "G" = 14.882091
    "H" = 44.646273
    "I" = 0.996401763178
    "J" = 0.996401763178
    "K" = 2
    "L" = 0.6666666666666666
    "M" = 0.2838346666666667
    "N" = 0.851504
    "O" = 0.7335492763611302
    "P" = 0.7335492763611302
    "Q" = 4.224053338288001
    "R" = 38.016480044592
    "S" = 0.7309097923441762
    "T" = 0.7309097923441762
}

$data[3] = @{
    "A" = "FAPs"
    "B" = "Tac1"
    "C" = "Tacr1"
    "D" = "FAPs"
    "E" = 3
    "F" = 1
    "G" = 14.882091
    "H" = 44.646273
    "I" = 0.996401763178
    "J" = 0.996401763178
    "K" = 1
    "L" = 0.3333333333333333
    "M" = 0.057966
    "N" = 0.173898
    "O" = 0.1498087525844245
    "P" = 0.1498087525844245
    "Q" = 0.862655286906
    "R" = 7.763897582154
    "S" = 0.1492697052146174
    "T" = 0.1492697052146174
}

$data[4] = @{
    "A" = "FAPs"
    "B" = "Tac1"
    "C" = "Tacr1"
    "D" = "MuSCs"
    "E" = 3
    "F" = 1
    "G" = 14.882091
    "H" = 44.646273
    "I" = 0.996401763178
    "J" = 0.996401763178
    "K" = 1
    "L" = 0.3333333333333333
    "M" = 0.04513266666666666
    "N" = 0.135398
    "O" = 0.1166419710544452
    "P" = 0.1166419710544452
    "Q" = 0.6716684524059999
    "R" = 6.045016071654
    "S" = 0.1162222656192064
    "T" = 0.1162222656192064
}

$data[5] = @{
    "A" = "MuSCs"
    "B" = "Tac1"
    "C" = "Tacr1"
    "D" = "ECs"
    "E" = 1
    "F" = 0.3333333333333333
    "G" = 0.037615
    "H" = 0.112845
    "I" = 0.002518439937098924
    "J" = 0.002518439937098924
    "K" = 2
    "L" = 0.6666666666666666
    "M" = 0.2838346666666667
    "N" = 0.851504
    "O" = 0.7335492763611302
    "P" = 0.7335492763611302
    "Q" = 0.01067644098666667
    "R" = 0.09608796888
    "S" = 0.001847399793417886
    "T" = 0.001847399793417886
}

$data[6] = @{
    "A" = "MuSCs"
    "B" = "Tac1"
    "C" = "Tacr1"
    "D" = "FAPs"
    "E" = 1
    "F" = 0.3333333333333333
    "G" = 0.037615
    "H" = 0.112845
    "I" = 0.002518439937098924
    "J" = 0.002518439937098924
    "K" = 1
    "L" = 0.3333333333333333
    "M" = 0.057966
    "N" = 0.173898
    "O" = 0.1498087525844245
    "P" = 0.1498087525844245
    "Q" = 0.00218039109
    "R" = 0.01962351981
    "S" = 0.0003772843454355864
    "T" = 0.0003772843454355864
}

$data[7] = @{
    "A" = "MuSCs"
    "B" = "Tac1"
    "C" = "Tacr1"
    "D" = "MuSCs"
    "E" = 1
    "F" = 0.3333333333333333
    "G" = 0.037615
    "H" = 0.112845
    "I" = 0.002518439937098924
    "J" = 0.002518439937098924
    "K" = 1
    "L" = 0.3333333333333333
    "M" = 0.04513266666666666
    "N" = 0.135398
    "O" = 0.1166419710544452
    "P" = 0.1166419710544452
    "Q" = 0.001697665256666667
    "R" = 0.01527898731
    "S" = 0.0002937557982454515
    "T" = 0.0002937557982454515
}

$data[8] = @{
    "A" = "Resolving-Mac"
    "B" = "Tac1"
    "C" = "Tacr1"
    "D" = "ECs"
    "E" = 1
    "F" = 0.3333333333333333
    "G" = 0.01612766666666667
    "H" = 0.048383
    "I" = 0.001079796884901035
    "J" = 0.001079796884901035
    "K" = 2
    "L" = 0.6666666666666666
    "M" = 0.2838346666666667
    "N" = 0.851504
    "O" = 0.7335492763611302
    "P" = 0.7335492763611302
    "Q" = 0.004577590892444445
    "R" = 0.041198318032
    "S" = 0.0007920842235361568
    "T" = 0.0007920842235361566
}

$data[9] = @{
    "A" = "Resolving-Mac"
    "B" = "Tac1"
    "C" = "Tacr1"
    "D" = "FAPs"
    "E" = 1
    "F" = 0.3333333333333333
    "G" = 0.01612766666666667
    "H" = 0.048383
    "I" = 0.001079796884901035
    "J" = 0.001079796884901035
    "K" = 1
    "L" = 0.3333333333333333
    "M" = 0.057966
    "N" = 0.173898
    "O" = 0.1498087525844245
    "P" = 0.1498087525844245
    "Q" = 0.0009348563260000001
    "R" = 0.008413706934
    "S" = 0.0001617630243715715
    "T" = 0.0001617630243715714
}

$data[10] = @{
    "A" = "Resolving-Mac"
    "B" = "Tac1"
    "C" = "Tacr1"
    "D" = "MuSCs"
    "E" = 1
    "F" = 0.3333333333333333
    "G" = 0.01612766666666667
    "H" = 0.048383
    "I" = 0.001079796884901035
    "J" = 0.001079796884901035
    "K" = 1
    "L" = 0.3333333333333333
    "M" = 0.04513266666666666
    "N" = 0.135398
    "O" = 0.1166419710544452
    "P" = 0.1166419710544452
    "Q" = 0.0007278846037777778
    "R" = 0.006550961434
    "S" = 0.0001259496369933066
    "T" = 0.0001259496369933066
}

foreach ($r in 2..10) {
    $row = $data[$r]
    foreach ($c in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")) {
        $ws.Range("$c$r").Value = $row[$c]
    }
}

